$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I20").Value = -0.4171720644817287
$ws.Range("J20").Value = 0.3357195371335476
$ws.Range("K20").Value = 0.2372264367509181
$ws.Range("L20").Value = 2.175867269640444
